$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(11, 1).Value = "Solitario"
$ws.Cells.Item(11, 2).Value = "Anna"
$ws.Cells.Item(11, 3).Value = 502
$ws.Cells.Item(11, 4).Value = 5
$ws.Cells.Item(11, 5).Value = "2025-11-08 20:50:26"
